$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Re-point the autofilter on column A ("Table") from "Products" to "IPG".
# This single call both rewrites the stored filter criterion and toggles
# row visibility: the "IPG" rows (24-26) become visible and the "Products"
# rows (40-74) become hidden, matching the new filter.
[void]$ws.Range("A1:G84").AutoFilter(1, "IPG", [Microsoft.Office.Interop.Excel.XlAutoFilterOperator]::xlFilterValues)

# Correct the "Inventory Posting Group" -> Code field's Fieldname value
# (was mistakenly "ItemID"). Done after the AutoFilter call above so the
# row is already unhidden when the cell is written.
$ws.Range("B24").Value = "Code"

# Widen column F so the "Inventory Posting Group" navtable label fits.
$ws.Columns.Item(6).ColumnWidth = 22

# Move the active selection to B25.
[void]$ws.Range("B25").Select()
